$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values like "58.462.22"
# or "6.95" are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.462.22"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "2.623.30"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "533.41"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").Value = "142.81"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "6.95"
$ws.Range("E9").Value = "  +7.04%  "

$ws.Range("E10").Value = "  -2.03%  "

$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").Value = "3.091.93"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").Value = "58.413.84"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("D15").Value = "20.72"
$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "2.619.83"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").Value = "4.39"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").Value = "334.34"
$ws.Range("E19").Value = "  -2.14%  "

$ws.Range("D20").Value = "10.11"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "66.32"
$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("D24").Value = "0.415"
$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "7.09"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").Value = "0.0$([char]0x2083)0736"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").Value = "5.86"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").Value = "18.76"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").Value = "150.28"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("D34").Value = "3.89"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "0.853"
$ws.Range("E35").Value = "  +2.11%  "

$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "1.41"
$ws.Range("E37").Value = "  -3.95%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "0.807"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "280.50"
$ws.Range("E40").Value = "  +2.58%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.595"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "10.68"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "18.99"
$ws.Range("E44").Value = "  +2.59%  "

$ws.Range("D45").Value = "0.0529"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0937"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0224"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.945.17"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.44"
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "17.84"
$ws.Range("E50").Value = "  -4.41%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "112.57"
$ws.Range("E51").Value = "  +1.33%  "

# Restore default styling for the Price/Volume columns (drop the temporary
# text number-format) while keeping the values stored as text.
$ws.Range("D2:E51").Style = "Normal"

